$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36 ("Iowa") now has successful data instead of an error.
$ws.Range("B36").Value = 44026
$ws.Range("B36").NumberFormat = $ws.Range("B37").NumberFormat
$ws.Range("C36").Value = 35865
$ws.Range("D36").Value = 757
$ws.Range("E36").Value = 3078
$ws.Range("F36").Value = 36
$ws.Range("G36").Value = 8.58
$ws.Range("H36").Value = 4.76
$ws.Range("I36").Value = $true
$ws.Range("J36").Value = $true
$ws.Range("O36").Value = "Success!"
